$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.053.25"
$ws.Range("E2").Value = "  +0.23%  "

# Row 3
$ws.Range("D3").Value = "1.823.38"
$ws.Range("E3").Value = "  +0.01%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.27%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.17"
$ws.Range("E5").Value = "  +0.35%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  +0.27%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4678"
$ws.Range("E7").Value = "  -0.19%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3644"
$ws.Range("E8").Value = "  -0.73%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07381"
$ws.Range("E9").Value = "  +0.26%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8776"
$ws.Range("E10").Value = "  +0.31%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.21"
$ws.Range("E11").Value = "  -0.41%  "

# Row 12
$ws.Range("D12").Value = "1.878.27"
$ws.Range("E12").Value = "  +1.87%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07570"
$ws.Range("E13").Value = "  +4.78%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.363"
$ws.Range("E14").Value = "  -1.25%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.71"
$ws.Range("E15").Value = "  +0.98%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.519"
$ws.Range("E16").Value = "  -0.08%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008714"
$ws.Range("E18").Value = "  -0.50%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("E19").Value = "  +0.27%  "

# Row 20
$ws.Range("D20").Value = "27.480.08"
$ws.Range("E20").Value = "  +1.77%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.59"
$ws.Range("E21").Value = "  -0.68%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.230"
$ws.Range("E22").Value = "  -1.20%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.60"
$ws.Range("E23").Value = "  -0.13%  "

# Row 24
$ws.Range("D24").Value = "2.080.36"
$ws.Range("E24").Value = "  +0.87%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.879"
$ws.Range("E25").Value = "  -0.63%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.30"
$ws.Range("E26").Value = "  +0.24%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.46"
$ws.Range("E27").Value = "  +0.51%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.123"
$ws.Range("E28").Value = "  -1.08%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.159"
$ws.Range("E29").Value = "  -1.63%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.39"
$ws.Range("E30").Value = "  -0.53%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08913"
$ws.Range("E31").Value = "  +0.59%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7427"
$ws.Range("E32").Value = "  -1.71%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.160"
$ws.Range("E33").Value = "  -0.32%  "

# Row 34
$ws.Range("E34").Value = "  -0.07%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.941"
$ws.Range("E35").Value = "  +0.21%  "

# Row 36
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.710"
$ws.Range("E36").Value = "  +13.78%  "

# Row 37
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.006"
$ws.Range("E37").Value = "  +0.30%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.087"
$ws.Range("E38").Value = "  -0.96%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05292"
$ws.Range("E39").Value = "  -0.39%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01927"
$ws.Range("E40").Value = "  -1.22%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.307"
$ws.Range("E41").Value = "  +1.52%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.929"
$ws.Range("E42").Value = "  -1.69%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5246"
$ws.Range("E43").Value = "  -1.16%  "

# Row 44
$ws.Range("E44").Value = "  -0.93%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.357"
$ws.Range("E45").Value = "  -1.53%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4902"
$ws.Range("E46").Value = "  +0.04%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.34"
$ws.Range("E47").Value = "  -1.05%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("E48").Value = "  +0.30%  "

# Row 49
$ws.Range("E49").Value = "  +1.09%  "

# Row 50
$ws.Range("E50").Value = "  -0.95%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06264"
$ws.Range("E51").Value = "  -0.53%  "

